$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = '43+20=63'
$t.Cell(1, 2).Range.Text = '46+29=75'
$t.Cell(1, 3).Range.Text = '55+28=83'
$t.Cell(1, 4).Range.Text = '56-30=26'
$t.Cell(1, 5).Range.Text = '21+51=72'
$t.Cell(2, 1).Range.Text = '22+8=30'
$t.Cell(2, 2).Range.Text = '99-95=4'
$t.Cell(2, 3).Range.Text = '41+16=57'
$t.Cell(2, 4).Range.Text = '16+56=72'
$t.Cell(2, 5).Range.Text = '65-49=16'
$t.Cell(3, 1).Range.Text = '44+4=48'
$t.Cell(3, 2).Range.Text = '50+2=52'
$t.Cell(3, 3).Range.Text = '24+22=46'
$t.Cell(3, 4).Range.Text = '96-35=61'
$t.Cell(3, 5).Range.Text = '52+31=83'
$t.Cell(4, 1).Range.Text = '17+74=91'
$t.Cell(4, 2).Range.Text = '33+30=63'
$t.Cell(4, 3).Range.Text = '5+83=88'
$t.Cell(4, 4).Range.Text = '51+18=69'
$t.Cell(4, 5).Range.Text = '64+10=74'
$t.Cell(5, 1).Range.Text = '55+2=57'
$t.Cell(5, 2).Range.Text = '13+53=66'
$t.Cell(5, 3).Range.Text = '98-53=45'
$t.Cell(5, 4).Range.Text = '62-16=46'
$t.Cell(5, 5).Range.Text = '83-4=79'
$t.Cell(6, 1).Range.Text = '51+1=52'
$t.Cell(6, 2).Range.Text = '81-59=22'
$t.Cell(6, 3).Range.Text = '6+37=43'
$t.Cell(6, 4).Range.Text = '18+24=42'
$t.Cell(6, 5).Range.Text = '19-0=19'
$t.Cell(7, 1).Range.Text = '20-15=5'
$t.Cell(7, 2).Range.Text = '74-31=43'
$t.Cell(7, 3).Range.Text = '67-1=66'
$t.Cell(7, 4).Range.Text = '44-20=24'
$t.Cell(7, 5).Range.Text = '36-34=2'
$t.Cell(8, 1).Range.Text = '66+23=89'
$t.Cell(8, 2).Range.Text = '76+0=76'
$t.Cell(8, 3).Range.Text = '20+11=31'
$t.Cell(8, 4).Range.Text = '98-45=53'
$t.Cell(8, 5).Range.Text = '12+26=38'
$t.Cell(9, 1).Range.Text = '20+12=32'
$t.Cell(9, 2).Range.Text = '58-27=31'
$t.Cell(9, 3).Range.Text = '77+10=87'
$t.Cell(9, 4).Range.Text = '12+1=13'
$t.Cell(9, 5).Range.Text = '81-43=38'
$t.Cell(10, 1).Range.Text = '17+12=29'
$t.Cell(10, 2).Range.Text = '94-88=6'
$t.Cell(10, 3).Range.Text = '25+63=88'
$t.Cell(10, 4).Range.Text = '1+23=24'
$t.Cell(10, 5).Range.Text = '18-1=17'
$t.Cell(11, 1).Range.Text = '80-78=2'
$t.Cell(11, 2).Range.Text = '47-44=3'
$t.Cell(11, 3).Range.Text = '17+62=79'
$t.Cell(11, 4).Range.Text = '47+20=67'
$t.Cell(11, 5).Range.Text = '13+81=94'
$t.Cell(12, 1).Range.Text = '15+57=72'
$t.Cell(12, 2).Range.Text = '4+92=96'
$t.Cell(12, 3).Range.Text = '74-17=57'
$t.Cell(12, 4).Range.Text = '30-11=19'
$t.Cell(12, 5).Range.Text = '39-34=5'
$t.Cell(13, 1).Range.Text = '75-2=73'
$t.Cell(13, 2).Range.Text = '40-15=25'
$t.Cell(13, 3).Range.Text = '19+44=63'
$t.Cell(13, 4).Range.Text = '10-10=0'
$t.Cell(13, 5).Range.Text = '14+0=14'
$t.Cell(14, 1).Range.Text = '25+73=98'
$t.Cell(14, 2).Range.Text = '72-0=72'
$t.Cell(14, 3).Range.Text = '34+4=38'
$t.Cell(14, 4).Range.Text = '33+31=64'
$t.Cell(14, 5).Range.Text = '27-8=19'
$t.Cell(15, 1).Range.Text = '28+5=33'
$t.Cell(15, 2).Range.Text = '50-18=32'
$t.Cell(15, 3).Range.Text = '16+61=77'
$t.Cell(15, 4).Range.Text = '3+94=97'
$t.Cell(15, 5).Range.Text = '88-24=64'
$t.Cell(16, 1).Range.Text = '70-55=15'
$t.Cell(16, 2).Range.Text = '97-57=40'
$t.Cell(16, 3).Range.Text = '59+31=90'
$t.Cell(16, 4).Range.Text = '21+22=43'
$t.Cell(16, 5).Range.Text = '3+61=64'
$t.Cell(17, 1).Range.Text = '81-17=64'
$t.Cell(17, 2).Range.Text = '1+31=32'
$t.Cell(17, 3).Range.Text = '36-21=15'
$t.Cell(17, 4).Range.Text = '1+34=35'
$t.Cell(17, 5).Range.Text = '86-42=44'
$t.Cell(18, 1).Range.Text = '87-47=40'
$t.Cell(18, 2).Range.Text = '5+76=81'
$t.Cell(18, 3).Range.Text = '30+51=81'
$t.Cell(18, 4).Range.Text = '2+60=62'
$t.Cell(18, 5).Range.Text = '39-4=35'
$t.Cell(19, 1).Range.Text = '26+72=98'
$t.Cell(19, 2).Range.Text = '17+69=86'
$t.Cell(19, 3).Range.Text = '0+93=93'
$t.Cell(19, 4).Range.Text = '16+15=31'
$t.Cell(19, 5).Range.Text = '1+60=61'
$t.Cell(20, 1).Range.Text = '24+68=92'
$t.Cell(20, 2).Range.Text = '97-94=3'
$t.Cell(20, 3).Range.Text = '76-9=67'
$t.Cell(20, 4).Range.Text = '47+50=97'
$t.Cell(20, 5).Range.Text = '41-27=14'
